# Split the sentence "A user may not be bothered with escaping:" into two
# runs (same rPr) so the first run keeps "A user ma" and the second run
# gets the reworded tail: "y not be bothered with escaping. Note it works
# only for Python 3:"

$d = $word.ActiveDocument

# Locate the original sentence inside the document body.
$sentence = $d.Content
$found = $sentence.Find.Execute("A user may not be bothered with escaping:", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the target sentence"
}

$sentenceStart = $sentence.Start
$sentenceEnd = $sentence.End

# Keep the leading "A user ma" (9 characters) in its own run - untouched -
# and turn everything after it into a distinct run carrying the new text.
$splitOffset = 9
$tail = $d.Range($sentenceStart + $splitOffset, $sentenceEnd)

$newTailXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr>' + `
    '<w:t>y not be bothered with escaping. Note it works only for Python 3:</w:t>' + `
    '</w:r></w:p>'

$tail.InsertXML($newTailXml)
